$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove any existing hyperlinks on the roster range so we can rebuild them
# cleanly in the new row order (avoids stale/duplicate relationship entries).
$ws.Range("K2:K18").Hyperlinks.Delete()

# Row 2: Franz Wagner
$ws.Cells.Item(2, 2).Value = 22
$ws.Cells.Item(2, 3).Value = "Franz Wagner"
$ws.Cells.Item(2, 4).Value = "SF"
$ws.Cells.Item(2, 5).Value = "6-9"
$ws.Cells.Item(2, 6).Value = 225
$ws.Cells.Item(2, 7).Value = "August 27, 2001"
$ws.Cells.Item(2, 8).Value = "de"
$ws.Cells.Item(2, 9).Value = "1"
$ws.Cells.Item(2, 10).Value = "Michigan"
$ws.Cells.Item(2, 11).Value = "https://www.basketball-reference.com/players/w/wagnefr01.html"

# Row 3: Bol Bol
$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(3, 3).Value = "Bol Bol"
$ws.Cells.Item(3, 4).Value = "PF"
$ws.Cells.Item(3, 5).Value = "7-2"
$ws.Cells.Item(3, 6).Value = 220
$ws.Cells.Item(3, 7).Value = "November 16, 1999"
$ws.Cells.Item(3, 8).Value = "sd"
$ws.Cells.Item(3, 9).Value = "3"
$ws.Cells.Item(3, 10).Value = "Oregon"
$ws.Cells.Item(3, 11).Value = "https://www.basketball-reference.com/players/b/bolbo01.html"

# Row 4: Paolo Banchero
$ws.Cells.Item(4, 2).Value = 5
$ws.Cells.Item(4, 3).Value = "Paolo Banchero"
$ws.Cells.Item(4, 4).Value = "PF"
$ws.Cells.Item(4, 5).Value = "6-10"
$ws.Cells.Item(4, 6).Value = 250
$ws.Cells.Item(4, 7).Value = "November 12, 2002"
$ws.Cells.Item(4, 8).Value = "us"
$ws.Cells.Item(4, 9).Value = "R"
$ws.Cells.Item(4, 10).Value = "Duke"
$ws.Cells.Item(4, 11).Value = "https://www.basketball-reference.com/players/b/banchpa01.html"

# Row 5: Terrence Ross
$ws.Cells.Item(5, 2).Value = 31
$ws.Cells.Item(5, 3).Value = "Terrence Ross"
$ws.Cells.Item(5, 4).Value = "SG"
$ws.Cells.Item(5, 5).Value = "6-6"
$ws.Cells.Item(5, 6).Value = 206
$ws.Cells.Item(5, 7).Value = "February 5, 1991"
$ws.Cells.Item(5, 8).Value = "us"
$ws.Cells.Item(5, 9).Value = "10"
$ws.Cells.Item(5, 10).Value = "Washington"
$ws.Cells.Item(5, 11).Value = "https://www.basketball-reference.com/players/r/rosste01.html"

# Row 6: Cole Anthony
$ws.Cells.Item(6, 2).Value = 50
$ws.Cells.Item(6, 3).Value = "Cole Anthony"
$ws.Cells.Item(6, 4).Value = "PG"
$ws.Cells.Item(6, 5).Value = "6-2"
$ws.Cells.Item(6, 6).Value = 185
$ws.Cells.Item(6, 7).Value = "May 15, 2000"
$ws.Cells.Item(6, 8).Value = "us"
$ws.Cells.Item(6, 9).Value = "2"
$ws.Cells.Item(6, 10).Value = "UNC"
$ws.Cells.Item(6, 11).Value = "https://www.basketball-reference.com/players/a/anthoco01.html"

# Row 7: Wendell Carter Jr.
$ws.Cells.Item(7, 2).Value = 34
$ws.Cells.Item(7, 3).Value = "Wendell Carter Jr."
$ws.Cells.Item(7, 4).Value = "C"
$ws.Cells.Item(7, 5).Value = "6-10"
$ws.Cells.Item(7, 6).Value = 270
$ws.Cells.Item(7, 7).Value = "April 16, 1999"
$ws.Cells.Item(7, 8).Value = "us"
$ws.Cells.Item(7, 9).Value = "4"
$ws.Cells.Item(7, 10).Value = "Duke"
$ws.Cells.Item(7, 11).Value = "https://www.basketball-reference.com/players/c/cartewe01.html"

# Row 8: Markelle Fultz
$ws.Cells.Item(8, 2).Value = 20
$ws.Cells.Item(8, 3).Value = "Markelle Fultz"
$ws.Cells.Item(8, 4).Value = "PG"
$ws.Cells.Item(8, 5).Value = "6-3"
$ws.Cells.Item(8, 6).Value = 209
$ws.Cells.Item(8, 7).Value = "May 29, 1998"
$ws.Cells.Item(8, 8).Value = "us"
$ws.Cells.Item(8, 9).Value = "5"
$ws.Cells.Item(8, 10).Value = "Washington"
$ws.Cells.Item(8, 11).Value = "https://www.basketball-reference.com/players/f/fultzma01.html"

# Row 9: Moritz Wagner
$ws.Cells.Item(9, 2).Value = 21
$ws.Cells.Item(9, 3).Value = "Moritz Wagner"
$ws.Cells.Item(9, 4).Value = "C"
$ws.Cells.Item(9, 5).Value = "6-11"
$ws.Cells.Item(9, 6).Value = 245
$ws.Cells.Item(9, 7).Value = "April 26, 1997"
$ws.Cells.Item(9, 8).Value = "de"
$ws.Cells.Item(9, 9).Value = "4"
$ws.Cells.Item(9, 10).Value = "Michigan"
$ws.Cells.Item(9, 11).Value = "https://www.basketball-reference.com/players/w/wagnemo01.html"

# Row 10: Jalen Suggs
$ws.Cells.Item(10, 2).Value = 4
$ws.Cells.Item(10, 3).Value = "Jalen Suggs"
$ws.Cells.Item(10, 4).Value = "PG"
$ws.Cells.Item(10, 5).Value = "6-4"
$ws.Cells.Item(10, 6).Value = 205
$ws.Cells.Item(10, 7).Value = "June 3, 2001"
$ws.Cells.Item(10, 8).Value = "us"
$ws.Cells.Item(10, 9).Value = "1"
$ws.Cells.Item(10, 10).Value = "Gonzaga"
$ws.Cells.Item(10, 11).Value = "https://www.basketball-reference.com/players/s/suggsja01.html"

# Row 11: Caleb Houstan
$ws.Cells.Item(11, 2).Value = 2
$ws.Cells.Item(11, 3).Value = "Caleb Houstan"
$ws.Cells.Item(11, 4).Value = "SF"
$ws.Cells.Item(11, 5).Value = "6-8"
$ws.Cells.Item(11, 6).Value = 205
$ws.Cells.Item(11, 7).Value = "January 9, 2003"
$ws.Cells.Item(11, 8).Value = "ca"
$ws.Cells.Item(11, 9).Value = "R"
$ws.Cells.Item(11, 10).Value = "Michigan"
$ws.Cells.Item(11, 11).Value = "https://www.basketball-reference.com/players/h/houstca01.html"

# Row 12: Gary Harris
$ws.Cells.Item(12, 2).Value = 14
$ws.Cells.Item(12, 3).Value = "Gary Harris"
$ws.Cells.Item(12, 4).Value = "SG"
$ws.Cells.Item(12, 5).Value = "6-4"
$ws.Cells.Item(12, 6).Value = 210
$ws.Cells.Item(12, 7).Value = "September 14, 1994"
$ws.Cells.Item(12, 8).Value = "us"
$ws.Cells.Item(12, 9).Value = "8"
$ws.Cells.Item(12, 10).Value = "Michigan State"
$ws.Cells.Item(12, 11).Value = "https://www.basketball-reference.com/players/h/harriga01.html"

# Row 13: Admiral Schofield (TW)
$ws.Cells.Item(13, 2).Value = 25
$ws.Cells.Item(13, 3).Value = "Admiral Schofield (TW)"
$ws.Cells.Item(13, 4).Value = "PF"
$ws.Cells.Item(13, 5).Value = "6-5"
$ws.Cells.Item(13, 6).Value = 241
$ws.Cells.Item(13, 7).Value = "March 30, 1997"
$ws.Cells.Item(13, 8).Value = "gb"
$ws.Cells.Item(13, 9).Value = "2"
$ws.Cells.Item(13, 10).Value = "Tennessee"
$ws.Cells.Item(13, 11).Value = "https://www.basketball-reference.com/players/s/schofad01.html"

# Row 14: Kevon Harris (TW)
$ws.Cells.Item(14, 2).Value = 7
$ws.Cells.Item(14, 3).Value = "Kevon Harris (TW)"
$ws.Cells.Item(14, 4).Value = "SG"
$ws.Cells.Item(14, 5).Value = "6-6"
$ws.Cells.Item(14, 6).Value = 216
$ws.Cells.Item(14, 7).Value = "June 24, 1997"
$ws.Cells.Item(14, 8).Value = "us"
$ws.Cells.Item(14, 9).Value = "R"
$ws.Cells.Item(14, 10).Value = "Stephen F. Austin"
$ws.Cells.Item(14, 11).Value = "https://www.basketball-reference.com/players/h/harrike01.html"

# Row 15: R.J. Hampton
$ws.Cells.Item(15, 2).Value = 13
$ws.Cells.Item(15, 3).Value = "R.J. Hampton"
$ws.Cells.Item(15, 4).Value = "PG"
$ws.Cells.Item(15, 5).Value = "6-4"
$ws.Cells.Item(15, 6).Value = 175
$ws.Cells.Item(15, 7).Value = "February 7, 2001"
$ws.Cells.Item(15, 8).Value = "us"
$ws.Cells.Item(15, 9).Value = "2"
$ws.Cells.Item(15, 10).Value = ""
$ws.Cells.Item(15, 11).Value = "https://www.basketball-reference.com/players/h/hamptrj01.html"

# Row 16: Chuma Okeke
$ws.Cells.Item(16, 2).Value = 3
$ws.Cells.Item(16, 3).Value = "Chuma Okeke"
$ws.Cells.Item(16, 4).Value = "SF"
$ws.Cells.Item(16, 5).Value = "6-6"
$ws.Cells.Item(16, 6).Value = 229
$ws.Cells.Item(16, 7).Value = "August 18, 1998"
$ws.Cells.Item(16, 8).Value = "us"
$ws.Cells.Item(16, 9).Value = "2"
$ws.Cells.Item(16, 10).Value = "Auburn"
$ws.Cells.Item(16, 11).Value = "https://www.basketball-reference.com/players/o/okekech01.html"

# Row 17: Jonathan Isaac
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = "Jonathan Isaac"
$ws.Cells.Item(17, 4).Value = "PF"
$ws.Cells.Item(17, 5).Value = "6-11"
$ws.Cells.Item(17, 6).Value = 230
$ws.Cells.Item(17, 7).Value = "October 3, 1997"
$ws.Cells.Item(17, 8).Value = "us"
$ws.Cells.Item(17, 9).Value = "3"
$ws.Cells.Item(17, 10).Value = "Florida State"
$ws.Cells.Item(17, 11).Value = "https://www.basketball-reference.com/players/i/isaacjo01.html"

# Row 18: Patrick Beverley
$ws.Cells.Item(18, 2).Value = ""
$ws.Cells.Item(18, 3).Value = "Patrick Beverley"
$ws.Cells.Item(18, 4).Value = "PG"
$ws.Cells.Item(18, 5).Value = "6-1"
$ws.Cells.Item(18, 6).Value = 180
$ws.Cells.Item(18, 7).Value = "July 12, 1988"
$ws.Cells.Item(18, 8).Value = "us"
$ws.Cells.Item(18, 9).Value = "10"
$ws.Cells.Item(18, 10).Value = "Arkansas"
$ws.Cells.Item(18, 11).Value = "https://www.basketball-reference.com/players/b/beverpa01.html"

# Rebuild hyperlinks in K2:K18 to match the new row order
$ws.Hyperlinks.Add($ws.Cells.Item(2, 11), "https://www.basketball-reference.com/players/w/wagnefr01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 11), "https://www.basketball-reference.com/players/b/bolbo01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 11), "https://www.basketball-reference.com/players/b/banchpa01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 11), "https://www.basketball-reference.com/players/r/rosste01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 11), "https://www.basketball-reference.com/players/a/anthoco01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 11), "https://www.basketball-reference.com/players/c/cartewe01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 11), "https://www.basketball-reference.com/players/f/fultzma01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9, 11), "https://www.basketball-reference.com/players/w/wagnemo01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 11), "https://www.basketball-reference.com/players/s/suggsja01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11, 11), "https://www.basketball-reference.com/players/h/houstca01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(12, 11), "https://www.basketball-reference.com/players/h/harriga01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(13, 11), "https://www.basketball-reference.com/players/s/schofad01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(14, 11), "https://www.basketball-reference.com/players/h/harrike01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(15, 11), "https://www.basketball-reference.com/players/h/hamptrj01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(16, 11), "https://www.basketball-reference.com/players/o/okekech01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(17, 11), "https://www.basketball-reference.com/players/i/isaacjo01.html") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(18, 11), "https://www.basketball-reference.com/players/b/beverpa01.html") | Out-Null
